# CTDC diagnosis 13 scripts
# Replaces the "StatQuery" and "CasesTab" query text on the startup sheet
# and adds a new "FilesTab" row describing the files query.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New query text for the StatQuery column (C2 / C3) ---
$statQuery = "MATCH (f:file)`nOPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)`nOPTIONAL MATCH (f)-[*]->(c:case)`nWITH f,a,ct,c`n   WHERE c.disease =  `"Adenocarcinoma of the gastroesophageal junction`"`nRETURN`n    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,`n    COUNT(DISTINCT c.case_id) AS Cases,`n    COUNT(DISTINCT f) AS Files"

# --- New CasesTab query text (B2) ---
$casesQuery = "MATCH (c:case)`n MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)`n MATCH (f:file)-[*]->(c)`n WHERE c.disease = `"Adenocarcinoma of the gastroesophageal junction`"`nRETURN DISTINCT`n    c.case_id AS ``Case ID``,`n     ct.clinical_trial_designation AS ``Trial Code``,`n     a.arm_id AS Arm,`n      a.arm_drug AS ``Arm Treatment``,`nc.disease AS Diagnosis,`n  c.gender AS Gender,`n    c.race AS Race,`n    c.ethnicity AS Ethnicity"

# --- New FilesTab query text (B3) ---
$filesQuery = "MATCH (f:file)`nOPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)`nOPTIONAL MATCH (f)-[*]->(c:case)`nOPTIONAL MATCH (f)-->(parent)`nWITH f,a,ct,c,parent`n WHERE c.disease = `"Adenocarcinoma of the gastroesophageal junction`"`nWITH`n    f, parent, c, a, ct,`n    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,`n    toInteger(floor(log(f.file_size)/log(1024))) as i,`n    2 as precision`nWITH`n    f, parent, c, a, ct,`n    f.file_size /(1024^i) AS value,`n    10^precision AS factor,`n    units[i] as unit`nWITH`n    f, parent, c, a, ct, unit,`n    round(factor * value)/factor AS size`nRETURN DISTINCT`n    f.file_name AS ``File Name``,`n    head(labels(parent)) as Association,`n    f.file_description AS Description,`n    f.file_format AS ``File Format``,`n    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,`n    ct.clinical_trial_designation AS ``Trial Code``,`n    a.arm_id AS Arm,`n    c.case_id AS ``Case ID``"

# Order of writes matters for shared-string ordering: StatQuery (C2) first,
# then the new FilesTab label (A3), then the CasesTab query text (B2),
# then the FilesTab query text (B3); C3 re-uses the StatQuery string.
$ws.Range("C2").Value() = $statQuery
$ws.Range("A3").Value() = "FilesTab"
$ws.Range("B2").Value() = $casesQuery
$ws.Range("B3").Value() = $filesQuery
$ws.Range("C3").Value() = $statQuery

$ws.Range("D3").Value() = $ws.Range("D2").Value()
$ws.Range("E3").Value() = $ws.Range("E2").Value()

# Wrap text on the new row's query cells (matches existing CasesTab formatting)
$ws.Range("B3").WrapText() = $true
$ws.Range("C3").WrapText() = $true

# Row heights to fit the new/expanded wrapped text
$ws.Rows.Item(2).RowHeight() = 195
$ws.Rows.Item(3).RowHeight() = 409.5

# Update the view: scrolled to the new row, with B3 selected
$ws.Application.ActiveWindow.ScrollRow() = 3
$ws.Range("B3").Select()
